$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.898.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.548.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.769.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.548.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.884.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0713"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.407.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.527"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.683.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("E51").Value = "  +0.25%  "
